$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.892.07"
$ws.Range("E2").Value = "  -4.28%  "
$ws.Range("D3").Value = "'1.636.33"
$ws.Range("E3").Value = "  -6.22%  "
$ws.Range("D4").Value = "'0.9979"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'231.77"
$ws.Range("E5").Value = "  -6.65%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "'0.4782"
$ws.Range("E7").Value = "  -5.36%  "
$ws.Range("D8").Value = "'39.19"
$ws.Range("E8").Value = "  -3.85%  "
$ws.Range("D9").Value = "'0.2582"
$ws.Range("E9").Value = "  -6.23%  "
$ws.Range("D10").Value = "'0.06093"
$ws.Range("E10").Value = "  -1.64%  "
$ws.Range("D11").Value = "'0.07033"
$ws.Range("E11").Value = "  -3.14%  "
$ws.Range("D12").Value = "'1.638.42"
$ws.Range("E12").Value = "  -6.10%  "
$ws.Range("D13").Value = "'14.42"
$ws.Range("E13").Value = "  -5.11%  "
$ws.Range("D14").Value = "'0.5939"
$ws.Range("E14").Value = "  -8.91%  "
$ws.Range("D15").Value = "'4.348"
$ws.Range("E15").Value = "  -7.40%  "
$ws.Range("D16").Value = "'73.29"
$ws.Range("E16").Value = "  -5.59%  "
$ws.Range("D17").Value = "'0.9997"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "'0.9991"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").Value = "'24.898.14"
$ws.Range("E19").Value = "  -4.32%  "
$ws.Range("D20").Value = "'0.000006571"
$ws.Range("E20").Value = "  -4.96%  "
$ws.Range("D21").Value = "'11.19"
$ws.Range("E21").Value = "  -6.08%  "
$ws.Range("D22").Value = "'1.848.35"
$ws.Range("E22").Value = "  -6.09%  "
$ws.Range("D23").Value = "'4.354"
$ws.Range("E23").Value = "  -2.89%  "
$ws.Range("D24").Value = "'8.535"
$ws.Range("E24").Value = "  -2.59%  "
$ws.Range("D25").Value = "'5.227"
$ws.Range("E25").Value = "  -2.62%  "
$ws.Range("D26").Value = "'133.25"
$ws.Range("E26").Value = "  -2.15%  "
$ws.Range("D27").Value = "'14.80"
$ws.Range("E27").Value = "  -3.21%  "
$ws.Range("D28").Value = "'1.383"
$ws.Range("E28").Value = "  -8.18%  "
$ws.Range("D29").Value = "'103.50"
$ws.Range("E29").Value = "  -2.18%  "
$ws.Range("D30").Value = "'1.631"
$ws.Range("E30").Value = "  -8.50%  "
$ws.Range("D31").Value = "'3.837"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").Value = "'0.07671"
$ws.Range("E32").Value = "  -6.48%  "
$ws.Range("D33").Value = "'3.522"
$ws.Range("E33").Value = "  -3.57%  "
$ws.Range("D34").Value = "'0.9988"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "'0.04281"
$ws.Range("E35").Value = "  -8.29%  "
$ws.Range("D36").Value = "'2.575"
$ws.Range("E36").Value = "  -2.93%  "
$ws.Range("D37").Value = "'0.9218"
$ws.Range("E37").Value = "  -7.51%  "
$ws.Range("D38").Value = "'0.5824"
$ws.Range("E38").Value = "  -4.49%  "
$ws.Range("D39").Value = "'2.550"
$ws.Range("E39").Value = "  -9.38%  "
$ws.Range("D40").Value = "'0.9989"
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.8336"
$ws.Range("E41").Value = "  +8.41%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.01515"
$ws.Range("E42").Value = "  -6.90%  "
$ws.Range("D43").Value = "'97.89"
$ws.Range("E43").Value = "  -2.78%  "
$ws.Range("D44").Value = "'1.726"
$ws.Range("E44").Value = "  -10.79%  "
$ws.Range("D45").Value = "'0.3679"
$ws.Range("E45").Value = "  -6.38%  "
$ws.Range("D46").Value = "'4.673"
$ws.Range("E46").Value = "  -6.72%  "
$ws.Range("D47").Value = "'0.05186"
$ws.Range("E47").Value = "  -2.60%  "
$ws.Range("D48").Value = "'6.027"
$ws.Range("E48").Value = "  -4.81%  "
$ws.Range("E49").Value = "  -7.43%  "
$ws.Range("D50").Value = "'29.15"
$ws.Range("E50").Value = "  -5.11%  "
$ws.Range("D51").Value = "'0.9998"
$ws.Range("E51").Value = "  -0.21%  "
